$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Cells.Item(1, 1).Value = "#_Agents"
$ws.Cells.Item(1, 2).Value = "Coverage"
$ws.Cells.Item(1, 3).Value = "Avg_Total_Rounds"
$ws.Cells.Item(1, 4).Value = "Avg_Expl_Cost"
$ws.Cells.Item(1, 5).Value = "Avg_Expl_Eff"
$ws.Cells.Item(1, 6).Value = "Avg_Round_Time"
$ws.Cells.Item(1, 7).Value = "Avg_Agent_Step_Time"
$ws.Cells.Item(1, 8).Value = "Avg_Experiment_Time"
$ws.Cells.Item(1, 9).Value = "Std_Total_Rounds"
$ws.Cells.Item(1, 10).Value = "Std_Expl_Cost"
$ws.Cells.Item(1, 11).Value = "Std_Expl_Eff"
$ws.Cells.Item(1, 12).Value = "Std_Round_Time"
$ws.Cells.Item(1, 13).Value = "Std_Agent_Step_Time"
$ws.Cells.Item(1, 14).Value = "Std_Experiment_Time"
$ws.Cells.Item(1, 15).Value = "Obs_Prob"

# Apply the same header style (bold, centered) used by the existing header
# cells (A1:J1) to the newly added header cells (K1:O1)
$ws.Range("A1").Copy()
$ws.Range("K1:O1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-13
# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 53.68
$ws.Cells.Item(2, 4).Value = 53.68
$ws.Cells.Item(2, 5).Value = 3.17812994
$ws.Cells.Item(2, 6).Value = 0.16001446
$ws.Cells.Item(2, 7).Value = 0.16001446
$ws.Cells.Item(2, 8).Value = 8.5570018
$ws.Cells.Item(2, 9).Value = 5.147678816164503
$ws.Cells.Item(2, 10).Value = 5.147678816164503
$ws.Cells.Item(2, 11).Value = 0.3145584481778171
$ws.Cells.Item(2, 12).Value = 0.01588193493321408
$ws.Cells.Item(2, 13).Value = 0.01588193493321408
$ws.Cells.Item(2, 14).Value = 0.9109293698347942
$ws.Cells.Item(2, 15).Value = 0.15

# Row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 86.63
$ws.Cells.Item(3, 4).Value = 86.63
$ws.Cells.Item(3, 5).Value = 1.98703528
$ws.Cells.Item(3, 6).Value = 0.11753604
$ws.Cells.Item(3, 7).Value = 0.11753604
$ws.Cells.Item(3, 8).Value = 10.07280072
$ws.Cells.Item(3, 9).Value = 11.95654707831956
$ws.Cells.Item(3, 10).Value = 11.95654707831956
$ws.Cells.Item(3, 11).Value = 0.268514152512527
$ws.Cells.Item(3, 12).Value = 0.02131027097002833
$ws.Cells.Item(3, 13).Value = 0.02131027097002833
$ws.Cells.Item(3, 14).Value = 1.760146952548398
$ws.Cells.Item(3, 15).Value = 0.85

# Row 4
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 31.004
$ws.Cells.Item(4, 4).Value = 61.97
$ws.Cells.Item(4, 5).Value = 2.78322638
$ws.Cells.Item(4, 6).Value = 0.2850598
$ws.Cells.Item(4, 7).Value = 0.14252968
$ws.Cells.Item(4, 8).Value = 4.396953860000001
$ws.Cells.Item(4, 9).Value = 4.475717606841654
$ws.Cells.Item(4, 10).Value = 8.919204247510965
$ws.Cells.Item(4, 11).Value = 0.3964916854618472
$ws.Cells.Item(4, 12).Value = 0.04660330493453917
$ws.Cells.Item(4, 13).Value = 0.0233015466708136
$ws.Cells.Item(4, 14).Value = 0.8670662738777775
$ws.Cells.Item(4, 15).Value = 0.15

# Row 5
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = 47.568
$ws.Cells.Item(5, 4).Value = 93.294
$ws.Cells.Item(5, 5).Value = 1.85185002
$ws.Cells.Item(5, 6).Value = 0.18283646
$ws.Cells.Item(5, 7).Value = 0.09141832
$ws.Cells.Item(5, 8).Value = 4.28448244
$ws.Cells.Item(5, 9).Value = 7.765951068091015
$ws.Cells.Item(5, 10).Value = 14.16137893379784
$ws.Cells.Item(5, 11).Value = 0.2725678101817861
$ws.Cells.Item(5, 12).Value = 0.03187944356968964
$ws.Cells.Item(5, 13).Value = 0.01593986758837946
$ws.Cells.Item(5, 14).Value = 0.7113420070683908
$ws.Cells.Item(5, 15).Value = 0.85

# Row 6
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = 18.638
$ws.Cells.Item(6, 4).Value = 74.416
$ws.Cells.Item(6, 5).Value = 2.32578502
$ws.Cells.Item(6, 6).Value = 0.36064252
$ws.Cells.Item(6, 7).Value = 0.09016050000000002
$ws.Cells.Item(6, 8).Value = 1.68161806
$ws.Cells.Item(6, 9).Value = 2.923705201317058
$ws.Cells.Item(6, 10).Value = 11.70114428400531
$ws.Cells.Item(6, 11).Value = 0.3594154407059264
$ws.Cells.Item(6, 12).Value = 0.08021947645837707
$ws.Cells.Item(6, 13).Value = 0.02005467182663801
$ws.Cells.Item(6, 14).Value = 0.463798327492031
$ws.Cells.Item(6, 15).Value = 0.15

# Row 7
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = 26.634
$ws.Cells.Item(7, 4).Value = 98.578
$ws.Cells.Item(7, 5).Value = 1.74627016
$ws.Cells.Item(7, 6).Value = 0.27419382
$ws.Cells.Item(7, 7).Value = 0.06854842
$ws.Cells.Item(7, 8).Value = 1.80959682
$ws.Cells.Item(7, 9).Value = 4.336086633658186
$ws.Cells.Item(7, 10).Value = 13.43640539653221
$ws.Cells.Item(7, 11).Value = 0.2387359758249279
$ws.Cells.Item(7, 12).Value = 0.06370485136771897
$ws.Cells.Item(7, 13).Value = 0.01592628321667663
$ws.Cells.Item(7, 14).Value = 0.4464294974852552
$ws.Cells.Item(7, 15).Value = 0.85

# Row 8
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = 12.668
$ws.Cells.Item(8, 4).Value = 75.888
$ws.Cells.Item(8, 5).Value = 2.31683128
$ws.Cells.Item(8, 6).Value = 0.41165586
$ws.Cells.Item(8, 7).Value = 0.06860936000000001
$ws.Cells.Item(8, 8).Value = 0.8771295800000001
$ws.Cells.Item(8, 9).Value = 2.308339850570419
$ws.Cells.Item(8, 10).Value = 13.8082871344412
$ws.Cells.Item(8, 11).Value = 0.5191305164231101
$ws.Cells.Item(8, 12).Value = 0.1192346384902147
$ws.Cells.Item(8, 13).Value = 0.01987252630628497
$ws.Cells.Item(8, 14).Value = 0.3294023284148944
$ws.Cells.Item(8, 15).Value = 0.15

# Row 9
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 18.666
$ws.Cells.Item(9, 4).Value = 95.774
$ws.Cells.Item(9, 5).Value = 1.80873238
$ws.Cells.Item(9, 6).Value = 0.28884228
$ws.Cells.Item(9, 7).Value = 0.04814038
$ws.Cells.Item(9, 8).Value = 0.8934506799999999
$ws.Cells.Item(9, 9).Value = 3.704306652472376
$ws.Cells.Item(9, 10).Value = 14.76912923488428
$ws.Cells.Item(9, 11).Value = 0.2942889472751435
$ws.Cells.Item(9, 12).Value = 0.07352609914658213
$ws.Cells.Item(9, 13).Value = 0.01225404545918741
$ws.Cells.Item(9, 14).Value = 0.2778735680592697
$ws.Cells.Item(9, 15).Value = 0.85

# Row 10
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 9.528
$ws.Cells.Item(10, 4).Value = 76.026
$ws.Cells.Item(10, 5).Value = 2.34660458
$ws.Cells.Item(10, 6).Value = 0.36560598
$ws.Cells.Item(10, 7).Value = 0.04570056
$ws.Cells.Item(10, 8).Value = 0.44229886
$ws.Cells.Item(10, 9).Value = 2.134086373439135
$ws.Cells.Item(10, 10).Value = 17.03248470042291
$ws.Cells.Item(10, 11).Value = 0.5813009885408598
$ws.Cells.Item(10, 12).Value = 0.1120422344700859
$ws.Cells.Item(10, 13).Value = 0.01400517053411987
$ws.Cells.Item(10, 14).Value = 0.1910490882520623
$ws.Cells.Item(10, 15).Value = 0.15

# Row 11
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = 14.542
$ws.Cells.Item(11, 4).Value = 90.11
$ws.Cells.Item(11, 5).Value = 1.94034866
$ws.Cells.Item(11, 6).Value = 0.26271068
$ws.Cells.Item(11, 7).Value = 0.03283892
$ws.Cells.Item(11, 8).Value = 0.4759576599999999
$ws.Cells.Item(11, 9).Value = 3.733109724670957
$ws.Cells.Item(11, 10).Value = 16.53410576071815
$ws.Cells.Item(11, 11).Value = 0.3659283612363473
$ws.Cells.Item(11, 12).Value = 0.07516614324360295
$ws.Cells.Item(11, 13).Value = 0.009395605954706286
$ws.Cells.Item(11, 14).Value = 0.1918604916329145
$ws.Cells.Item(11, 15).Value = 0.85

# Row 12
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = 7.264
$ws.Cells.Item(12, 4).Value = 72.416
$ws.Cells.Item(12, 5).Value = 2.49796812
$ws.Cells.Item(12, 6).Value = 0.35798356
$ws.Cells.Item(12, 7).Value = 0.0357984
$ws.Cells.Item(12, 8).Value = 0.26875334
$ws.Cells.Item(12, 9).Value = 1.828932472229621
$ws.Cells.Item(12, 10).Value = 18.25248614575565
$ws.Cells.Item(12, 11).Value = 0.6916429980740318
$ws.Cells.Item(12, 12).Value = 0.1215264625450928
$ws.Cells.Item(12, 13).Value = 0.01215267791892809
$ws.Cells.Item(12, 14).Value = 0.1392939107174838
$ws.Cells.Item(12, 15).Value = 0.15

# Row 13
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = 11.936
$ws.Cells.Item(13, 4).Value = 84.272
$ws.Cells.Item(13, 5).Value = 2.08262574
$ws.Cells.Item(13, 6).Value = 0.23579416
$ws.Cells.Item(13, 7).Value = 0.02357946
$ws.Cells.Item(13, 8).Value = 0.2809738
$ws.Cells.Item(13, 9).Value = 3.304806725814685
$ws.Cells.Item(13, 10).Value = 16.52198438250901
$ws.Cells.Item(13, 11).Value = 0.4083749073524997
$ws.Cells.Item(13, 12).Value = 0.06549908709527966
$ws.Cells.Item(13, 13).Value = 0.006549929454591113
$ws.Cells.Item(13, 14).Value = 0.11199575969267
$ws.Cells.Item(13, 15).Value = 0.85
